$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add two new worksheets (Sheet3, Sheet4) after the existing sheets
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Sheet3"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "Sheet4"

# ---------------------------------------------------------------------
# 2. Sheet3 - quick scratch check of the two test fixtures
# ---------------------------------------------------------------------
$ws3.Range("A1").Value = "Actual"
$ws3.Range("A3").Value = "[[[0, 0], [0, 1], [0, 2], [0, 3]], [[0, 0], [0, 0], [0, 1], [0, 2]], [[0, 0], [0, 0], [0, 0], [0, 1]], [[0, 0], [0, 0], [0, 0], [0, 0]], [[0, 0], [0, 0], [0, 0], [0, 0]], [[0, 0], [0, 0], [0, 0], [0, 0]], [[0, 0], [0, 0], [0, 0], [0, 0]], [[0, 0], [0, 0], [0, 0], [0, 0]], [[0, 0], [0, 0], [0, 0], [0, 0]]]"
$ws3.Range("A5").Value = "[[[0, 0], [0, 1], [0, 2], [0, 3]], [[0, 0], [0, 0], [0, 1], [0, 2]], [[0, 0], [0, 0], [0, 0], [0, 0]], [[0, 0], [0, 0], [0, 0], [0, 0]], [[0, 0], [0, 0], [0, 0], [0, 0]], [[0, 0], [0, 0], [0, 0], [0, 0]], [[0, 0], [0, 0], [0, 0], [0, 0]], [[0, 0], [0, 0], [0, 0], [0, 0]]]"
$ws3.Range("M11").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Sheet4 - actual vs expected comparison with formulas
# ---------------------------------------------------------------------
$ws4.Range("A1").Value = "Actual"

$ws4.Range("A3").Value = "[[0, 0], [0, 1], [0, 2], [0, 3]]"
$ws4.Range("A4").Value = "[[0, 0], [0, 0], [0, 1], [0, 2]]"
$ws4.Range("A5").Value = "[[0, 0], [0, 0], [0, 0], [0, 1]]"
$ws4.Range("A6").Value = "[[0, 0], [0, 0], [0, 0], [0, 0]]"
$ws4.Range("A7").Value = "[[0, 0], [0, 0], [0, 0], [0, 0]]"
$ws4.Range("A8").Value = "[[0, 0], [0, 0], [0, 0], [0, 0]]"
$ws4.Range("A9").Value = "[[0, 0], [0, 0], [0, 0], [0, 0]]"
$ws4.Range("A10").Value = "[[0, 0], [0, 0], [0, 0], [0, 0]]"
$ws4.Range("A11").Value = "[[0, 0], [0, 0], [0, 0], [0, 0]]"

$ws4.Range("A13").Value = "Expected"

$ws4.Range("A15").Value = "[[0, 0], [0, 1], [0, 2], [0, 3]]"
$ws4.Range("A16").Value = "[[0, 0], [0, 0], [0, 1], [0, 2]]"
$ws4.Range("A17").Value = "[[0, 0], [0, 0], [0, 0], [0, 0]]"
$ws4.Range("A18").Value = "[[0, 0], [0, 0], [0, 0], [0, 0]]"
$ws4.Range("A19").Value = "[[0, 0], [0, 0], [0, 0], [0, 0]]"
$ws4.Range("A20").Value = "[[0, 0], [0, 0], [0, 0], [0, 0]]"
$ws4.Range("A21").Value = "[[0, 0], [0, 0], [0, 0], [0, 0]]"
$ws4.Range("A22").Value = "[[0, 0], [0, 0], [0, 0], [0, 0]]"
$ws4.Range("A23").Value = "[[0, 0], [0, 0], [0, 0], [0, 0]]"

$ws4.Range("B3").Formula = "=(A3=A15)"
$ws4.Range("B4:B11").Formula = "=(A4=A16)"

$ws4.Range("A17").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. Sheet2 - move the "visited" highlight from G12 to D10 and B12,
#    and add a new grey highlight colour for the second marker
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

# Clear G12 back to the plain bordered style by copying a known-plain
# neighbouring cell's formatting onto it.
$ws2.Range("H12").Copy($ws2.Range("G12")) | Out-Null

# Give D10 a new grey fill (theme "Background 2, Darker 25%") and copy
# that exact formatting onto B12 as well.
$d10 = $ws2.Range("D10")
$d10.Interior.Color = 7434614
$d10.Copy($ws2.Range("B12")) | Out-Null

$ws2.Range("G18").Select() | Out-Null

# ---------------------------------------------------------------------
# 5. Make Sheet2 the active tab
# ---------------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("G18").Select() | Out-Null
